$p = $ppt.ActivePresentation

$newDate = "1/25/2023"

# --- Update the cached "datetimeFigureOut" field text wherever it appears
# (slide master, all slide layouts, handout master, notes master) ---

$slideMaster = $p.SlideMaster
$slideMaster.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

$layouts = $slideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$handoutMaster = $p.HandoutMaster
for ($j = 1; $j -le $handoutMaster.Shapes.Count; $j++) {
    $shp = $handoutMaster.Shapes.Item($j)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

$notesMaster = $p.NotesMaster
for ($j = 1; $j -le $notesMaster.Shapes.Count; $j++) {
    $shp = $notesMaster.Shapes.Item($j)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# --- Slide 1: update the byline text from "Terence Parr" to "Mustafa Hajij" ---

$slide1 = $p.Slides.Item(1)
for ($j = 1; $j -le $slide1.Shapes.Count; $j++) {
    $shp = $slide1.Shapes.Item($j)
    if ($shp.Name -eq "Rectangle 4") {
        $shp.TextFrame.TextRange.Paragraphs(1).Text = "Mustafa Hajij"
    }
}
